$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2,3,4,5,7,8 form a cycle; each row's data (columns D,L,M,N,O,P,Q,S,T)
# gets replaced by the data currently held by the next row in the cycle
# (row 6 is left untouched). Capture all "before" values first, then write
# them in a second pass so the in-place shift doesn't clobber source data.

$cycle = @(2, 3, 4, 5, 7, 8)

$data = @{}
foreach ($r in $cycle) {
    $data[$r] = @{
        D = $ws.Cells.Item($r, 4).Value2
        L = $ws.Cells.Item($r, 12).Value2
        M = $ws.Cells.Item($r, 13).Value2
        N = $ws.Cells.Item($r, 14).Value2
        O = $ws.Cells.Item($r, 15).Value2
        P = $ws.Cells.Item($r, 16).Value2
        Q = $ws.Cells.Item($r, 17).Value2
        S = $ws.Cells.Item($r, 19).Value2
        T = $ws.Cells.Item($r, 20).Value2
    }
}

for ($i = 0; $i -lt $cycle.Length; $i++) {
    $destRow = $cycle[$i]
    $srcRow = $cycle[($i + 1) % $cycle.Length]
    $src = $data[$srcRow]

    $ws.Cells.Item($destRow, 4).Value = $src.D
    $ws.Cells.Item($destRow, 12).Value = $src.L
    $ws.Cells.Item($destRow, 13).Value = $src.M
    $ws.Cells.Item($destRow, 14).Value = $src.N
    $ws.Cells.Item($destRow, 15).Value = $src.O
    $ws.Cells.Item($destRow, 16).Value = $src.P
    $ws.Cells.Item($destRow, 17).Value = $src.Q
    $ws.Cells.Item($destRow, 19).Value = $src.S
    $ws.Cells.Item($destRow, 20).Value = $src.T
}
